# Slide 1, TextBox "TextBox 8": remove the middle run "복수키" (err="1")
# and change the trailing run text from " 한 번에 삭제" to "데이터 한 번에 삭제",
# so the combined text goes from "Redis 복수키 한 번에 삭제" to
# "Redis 데이터 한 번에 삭제".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# Locate and delete the "복수키" run (characters 7-9 of "Redis 복수키 한 번에 삭제").
$midRun = $tr.Characters(7, 3)
$midRun.Delete()

# After deletion the remaining tail run (" 한 번에 삭제") starts at character 7;
# replace its text so the run reads "데이터 한 번에 삭제".
$tailRun = $tr.Characters(7, $tr.Length - 6)
$tailRun.Text = "데이터 한 번에 삭제"
